$wb = $excel.ActiveWorkbook

# Column E width on "Design 1 Data" changes from 21 to 22 (stored OOXML width).
# Excel's ColumnWidth property is offset from the stored width by the default
# font padding (~0.8333 chars for Calibri 11); 21.166666666666668 round-trips to
# a saved width of exactly 22.
$wsDesign1 = $wb.Worksheets.Item("Design 1 Data")
$wsDesign1.Columns.Item(5).ColumnWidth = 21.166666666666668

# ---- Design 1 Data ----
$ws_Design_1_Data = $wb.Worksheets.Item("Design 1 Data")
$ws_Design_1_Data.Range("K2").Value = 2617230.057129697
$ws_Design_1_Data.Range("K3").Value = 1285604.607532217
$ws_Design_1_Data.Range("K4").Value = 1285179.607532217
$ws_Design_1_Data.Range("K5").Value = 2175291.087589347
$ws_Design_1_Data.Range("K6").Value = 0.07818974950738743
$ws_Design_1_Data.Range("K7").Value = 1.614820666231021
$ws_Design_1_Data.Range("H8").Value = 13.6023089546819
$ws_Design_1_Data.Range("K8").Value = 0.08478783657025808
$ws_Design_1_Data.Range("H9").Value = 4837280592.753655
$ws_Design_1_Data.Range("H10").Value = 0.001270043469758956
$ws_Design_1_Data.Range("K10").Value = 4999.962615982123
$ws_Design_1_Data.Range("K11").Value = 30867989.59613708
$ws_Design_1_Data.Range("E12").Value = 0.0180407065178009
$ws_Design_1_Data.Range("H12").Value = 46205736.73514064
$ws_Design_1_Data.Range("J13").Value = "total_fuel_L [L]"
$ws_Design_1_Data.Range("K13").Value = 54938.8341339537
$ws_Design_1_Data.Range("H14").Value = 64.8165
$ws_Design_1_Data.Range("J14").Value = "mission_fuel_L [L]"
$ws_Design_1_Data.Range("K14").Value = 54385.67286980832
$ws_Design_1_Data.Range("B15").Value = 61.73
$ws_Design_1_Data.Range("J15").Value = "reserve_fuel_L [L]"
$ws_Design_1_Data.Range("K15").Value = 553.1612641453846
$ws_Design_1_Data.Range("B16").Value = 61.73
$ws_Design_1_Data.Range("J16").Value = "max_fuel_L [L]"
$ws_Design_1_Data.Range("K16").Value = 60432.71754734909
$ws_Design_1_Data.Range("J17").Value = "MTOM [kg]"
$ws_Design_1_Data.Range("K17").Value = 266792.0547532821
$ws_Design_1_Data.Range("H18").Value = 134.0936683963558
$ws_Design_1_Data.Range("J18").Value = "S [m²]"
$ws_Design_1_Data.Range("K18").Value = 523.4200260289563
$ws_Design_1_Data.Range("H19").Value = 38.14450654420889
$ws_Design_1_Data.Range("J19").Value = "b [m]"
$ws_Design_1_Data.Range("K19").Value = 51.15760090294287
$ws_Design_1_Data.Range("J20").Value = "MAC [m]"
$ws_Design_1_Data.Range("K20").Value = 10.23152018058857
$ws_Design_1_Data.Range("H21").Value = 24.56463937825266
$ws_Design_1_Data.Range("J21").Value = "fuel_economy [L/ton/km]"
$ws_Design_1_Data.Range("K21").Value = 0.07958242297173541
$ws_Design_1_Data.Range("H22").Value = 8.398167308804327
$ws_Design_1_Data.Range("J22").Value = "OEW [N]"
$ws_Design_1_Data.Range("K22").Value = 1285604.607532217
$ws_Design_1_Data.Range("H23").Value = 2.519450192641298
$ws_Design_1_Data.Range("J23").Value = "total_fuel [N]"
$ws_Design_1_Data.Range("K23").Value = 441938.9695403504
$ws_Design_1_Data.Range("J24").Value = "max_fuel [N]"
$ws_Design_1_Data.Range("K24").Value = 486132.8664943855
$ws_Design_1_Data.Range("J25").Value = "mission_fuel [N]"
$ws_Design_1_Data.Range("K25").Value = 437489.2296993121
$ws_Design_1_Data.Range("J26").Value = "reserve_fuel [N]"
$ws_Design_1_Data.Range("K26").Value = 4449.739841038303
$ws_Design_1_Data.Range("K26").NumberFormat = "#,##0.000"
$ws_Design_1_Data.Range("H27").Value = 5.986385927814366
$ws_Design_1_Data.Range("J27").Value = "Mff"
$ws_Design_1_Data.Range("K27").Value = 0.8117510033563166
$ws_Design_1_Data.Range("K27").NumberFormat = "#,##0.000"
$ws_Design_1_Data.Range("J28").Value = "LD"
$ws_Design_1_Data.Range("K28").Value = 21.9652896084792
$ws_Design_1_Data.Range("K28").NumberFormat = "#,##0.000"
$ws_Design_1_Data.Range("J29").Value = "MTOW [N]"
$ws_Design_1_Data.Range("K29").Value = 2617230.057129697
$ws_Design_1_Data.Range("K29").NumberFormat = "#,##0.000"
$ws_Design_1_Data.Range("H30").Value = 18.35485255942574
$ws_Design_1_Data.Range("J30").Value = "total_fuel_L [L]"
$ws_Design_1_Data.Range("K30").Value = 54938.8341339537
$ws_Design_1_Data.Range("K30").NumberFormat = "#,##0.000"
$ws_Design_1_Data.Range("H31").Value = 66.86495379374394
$ws_Design_1_Data.Range("J31").Value = "mission_fuel_L [L]"
$ws_Design_1_Data.Range("K31").Value = 54385.67286980832
$ws_Design_1_Data.Range("K31").NumberFormat = "#,##0.000"
$ws_Design_1_Data.Range("H32").Value = 0.9363014827246962
$ws_Design_1_Data.Range("J32").Value = "reserve_fuel_L [L]"
$ws_Design_1_Data.Range("K32").Value = 553.1612641453846
$ws_Design_1_Data.Range("K32").NumberFormat = "#,##0.000"
$ws_Design_1_Data.Range("J33").Value = "max_fuel_L [L]"
$ws_Design_1_Data.Range("K33").Value = 60432.71754734907
$ws_Design_1_Data.Range("K33").NumberFormat = "#,##0.000"
$ws_Design_1_Data.Range("H35").Value = 60.44726267805435
$ws_Design_1_Data.Range("H36").Value = 37.15138198133811
$ws_Design_1_Data.Range("H38").Value = 9.522126549100339
$ws_Design_1_Data.Range("H39").Value = 7.468334548313991
$ws_Design_1_Data.Range("H40").Value = 5.227834183819794
$ws_Design_1_Data.Range("E42").Value = 14.61645740084082
$ws_Design_1_Data.Range("E43").Value = 5.846582960336328
$ws_Design_1_Data.Range("H43").Value = 6.41398143561084
$ws_Design_1_Data.Range("E44").Value = 10.96234305063062
$ws_Design_1_Data.Range("H45").Value = 67.3098365394519
$ws_Design_1_Data.Range("E46").Value = 29.09163903076659
$ws_Design_1_Data.Range("H46").Value = 0.9231497192690139
$ws_Design_1_Data.Range("E47").Value = 10.85793978348175
$ws_Design_1_Data.Range("E48").Value = 523.4200260289563
$ws_Design_1_Data.Range("E49").Value = 51.15760090294287
$ws_Design_1_Data.Range("H49").Value = 32.1253565723892

# ---- Design 2 Data ----
$ws_Design_2_Data = $wb.Worksheets.Item("Design 2 Data")
$ws_Design_2_Data.Range("J13").Value = "max_fuel_L [L]"
$ws_Design_2_Data.Range("K13").Value = 50532.6180386827
$ws_Design_2_Data.Range("J14").Value = "total_fuel_L [L]"
$ws_Design_2_Data.Range("K14").Value = 45938.74367152973
$ws_Design_2_Data.Range("J15").Value = "mission_fuel_L [L]"
$ws_Design_2_Data.Range("K15").Value = 45560.14901022315
$ws_Design_2_Data.Range("J16").Value = "reserve_fuel_L [L]"
$ws_Design_2_Data.Range("K16").Value = 378.5946613065801
$ws_Design_2_Data.Range("J17").Value = "MTOM [kg]"
$ws_Design_2_Data.Range("K17").Value = 255365.0821417169
$ws_Design_2_Data.Range("J18").Value = "S [m²]"
$ws_Design_2_Data.Range("K18").Value = 498.4583613440194
$ws_Design_2_Data.Range("J19").Value = "b [m]"
$ws_Design_2_Data.Range("K19").Value = 70.60158364682901
$ws_Design_2_Data.Range("J20").Value = "MAC [m]"
$ws_Design_2_Data.Range("K20").Value = 7.060158364682901
$ws_Design_2_Data.Range("J21").Value = "fuel_economy [L/ton/km]"
$ws_Design_2_Data.Range("K21").Value = 0.06666805535102044
$ws_Design_2_Data.Range("J22").Value = "OEW [N]"
$ws_Design_2_Data.Range("K22").Value = 1246016.632511913
$ws_Design_2_Data.Range("J23").Value = "total_fuel [N]"
$ws_Design_2_Data.Range("K23").Value = 369540.4418425194
$ws_Design_2_Data.Range("J24").Value = "max_fuel [N]"
$ws_Design_2_Data.Range("K24").Value = 406494.4860267714
$ws_Design_2_Data.Range("J25").Value = "mission_fuel [N]"
$ws_Design_2_Data.Range("K25").Value = 366494.9506680371
$ws_Design_2_Data.Range("J26").Value = "reserve_fuel [N]"
$ws_Design_2_Data.Range("K26").Value = 3045.491174482391
$ws_Design_2_Data.Range("K26").NumberFormat = "#,##0.000"
$ws_Design_2_Data.Range("J27").Value = "Mff"
$ws_Design_2_Data.Range("K27").Value = 0.8524866066467567
$ws_Design_2_Data.Range("K27").NumberFormat = "#,##0.000"
$ws_Design_2_Data.Range("J28").Value = "MTOW [N]"
$ws_Design_2_Data.Range("K28").Value = 2505131.455810243
$ws_Design_2_Data.Range("K28").NumberFormat = "#,##0.000"
$ws_Design_2_Data.Range("J29").Value = "LD"
$ws_Design_2_Data.Range("K29").Value = 31.49782606412045
$ws_Design_2_Data.Range("K29").NumberFormat = "#,##0.000"
$ws_Design_2_Data.Range("J30").Value = "total_fuel_L [L]"
$ws_Design_2_Data.Range("K30").Value = 45938.74367152973
$ws_Design_2_Data.Range("K30").NumberFormat = "#,##0.000"
$ws_Design_2_Data.Range("J31").Value = "mission_fuel_L [L]"
$ws_Design_2_Data.Range("K31").Value = 45560.14901022315
$ws_Design_2_Data.Range("K31").NumberFormat = "#,##0.000"
$ws_Design_2_Data.Range("J32").Value = "reserve_fuel_L [L]"
$ws_Design_2_Data.Range("K32").Value = 378.5946613065801
$ws_Design_2_Data.Range("K32").NumberFormat = "#,##0.000"
$ws_Design_2_Data.Range("J33").Value = "max_fuel_L [L]"
$ws_Design_2_Data.Range("K33").Value = 50532.6180386827
$ws_Design_2_Data.Range("K33").NumberFormat = "#,##0.000"

# ---- Design 3 Data ----
$ws_Design_3_Data = $wb.Worksheets.Item("Design 3 Data")
$ws_Design_3_Data.Range("K2").Value = 2538415.641737442
$ws_Design_3_Data.Range("K3").Value = 1257771.047870782
$ws_Design_3_Data.Range("K4").Value = 1257346.047870782
$ws_Design_3_Data.Range("K5").Value = 2147378.71351252
$ws_Design_3_Data.Range("K6").Value = 0.1568692251654692
$ws_Design_3_Data.Range("K7").Value = 1.655258149774672
$ws_Design_3_Data.Range("K8").Value = 0.1056332659449669
$ws_Design_3_Data.Range("K11").Value = 24030456.87387829
$ws_Design_3_Data.Range("E12").Value = 0.01888976835334421
$ws_Design_3_Data.Range("J13").Value = "max_fuel_L [L]"
$ws_Design_3_Data.Range("K13").Value = 53472.14403513278
$ws_Design_3_Data.Range("H14").Value = 16.81458221126043
$ws_Design_3_Data.Range("J14").Value = "total_fuel_L [L]"
$ws_Design_3_Data.Range("K14").Value = 48611.04003193889
$ws_Design_3_Data.Range("J15").Value = "mission_fuel_L [L]"
$ws_Design_3_Data.Range("K15").Value = 48180.18163096658
$ws_Design_3_Data.Range("J16").Value = "reserve_fuel_L [L]"
$ws_Design_3_Data.Range("K16").Value = 430.858400972309
$ws_Design_3_Data.Range("J17").Value = "MTOM [kg]"
$ws_Design_3_Data.Range("K17").Value = 258757.9655185976
$ws_Design_3_Data.Range("J18").Value = "S [m²]"
$ws_Design_3_Data.Range("K18").Value = 507.9658430584618
$ws_Design_3_Data.Range("H19").Value = 97.89060565452871
$ws_Design_3_Data.Range("J19").Value = "b [m]"
$ws_Design_3_Data.Range("K19").Value = 63.7473665688842
$ws_Design_3_Data.Range("H20").Value = 39.49258206671682
$ws_Design_3_Data.Range("J20").Value = "MAC [m]"
$ws_Design_3_Data.Range("K20").Value = 7.968420821110524
$ws_Design_3_Data.Range("H21").Value = 20.98827590454679
$ws_Design_3_Data.Range("J21").Value = "fuel_economy [L/ton/km]"
$ws_Design_3_Data.Range("K21").Value = 0.07050194271916771
$ws_Design_3_Data.Range("H22").Value = 7.175478941725396
$ws_Design_3_Data.Range("J22").Value = "OEW [N]"
$ws_Design_3_Data.Range("K22").Value = 1257771.047870782
$ws_Design_3_Data.Range("H23").Value = 2.152643682517619
$ws_Design_3_Data.Range("J23").Value = "total_fuel [N]"
$ws_Design_3_Data.Range("K23").Value = 391036.9282249228
$ws_Design_3_Data.Range("J24").Value = "max_fuel [N]"
$ws_Design_3_Data.Range("K24").Value = 430140.6210474151
$ws_Design_3_Data.Range("J25").Value = "mission_fuel [N]"
$ws_Design_3_Data.Range("K25").Value = 387571.0170758213
$ws_Design_3_Data.Range("J26").Value = "reserve_fuel [N]"
$ws_Design_3_Data.Range("K26").Value = 3465.911149101448
$ws_Design_3_Data.Range("K26").NumberFormat = "#,##0.000"
$ws_Design_3_Data.Range("H27").Value = 5.114828578973488
$ws_Design_3_Data.Range("J27").Value = "Mff"
$ws_Design_3_Data.Range("K27").Value = 0.8438028261990694
$ws_Design_3_Data.Range("K27").NumberFormat = "#,##0.000"
$ws_Design_3_Data.Range("J28").Value = "MTOW [N]"
$ws_Design_3_Data.Range("K28").Value = 2538415.641737442
$ws_Design_3_Data.Range("K28").NumberFormat = "#,##0.000"
$ws_Design_3_Data.Range("E29").Value = 0.08
$ws_Design_3_Data.Range("J29").Value = "LD"
$ws_Design_3_Data.Range("K29").Value = 27.83247462486353
$ws_Design_3_Data.Range("K29").NumberFormat = "#,##0.000"
$ws_Design_3_Data.Range("E30").Value = 0.08
$ws_Design_3_Data.Range("J30").Value = "total_fuel_L [L]"
$ws_Design_3_Data.Range("K30").Value = 48611.04003193889
$ws_Design_3_Data.Range("K30").NumberFormat = "#,##0.000"
$ws_Design_3_Data.Range("H31").Value = 67.60269214604048
$ws_Design_3_Data.Range("J31").Value = "mission_fuel_L [L]"
$ws_Design_3_Data.Range("K31").Value = 48180.18163096658
$ws_Design_3_Data.Range("K31").NumberFormat = "#,##0.000"
$ws_Design_3_Data.Range("H32").Value = 0.9505700135728937
$ws_Design_3_Data.Range("J32").Value = "reserve_fuel_L [L]"
$ws_Design_3_Data.Range("K32").Value = 430.858400972309
$ws_Design_3_Data.Range("K32").NumberFormat = "#,##0.000"
$ws_Design_3_Data.Range("J33").Value = "max_fuel_L [L]"
$ws_Design_3_Data.Range("K33").Value = 53472.14403513278
$ws_Design_3_Data.Range("K33").NumberFormat = "#,##0.000"
$ws_Design_3_Data.Range("H35").Value = 74.79953332310792
$ws_Design_3_Data.Range("H36").Value = 36.30688408377614
$ws_Design_3_Data.Range("H37").Value = 10.59241709831434
$ws_Design_3_Data.Range("H38").Value = 8.82701424859528
$ws_Design_3_Data.Range("H39").Value = 5.296208549157168
$ws_Design_3_Data.Range("E41").Value = 11.38345831587218
$ws_Design_3_Data.Range("E42").Value = 4.553383326348871
$ws_Design_3_Data.Range("E43").Value = 13.66014997904661
$ws_Design_3_Data.Range("H43").Value = 65.95115683917061
$ws_Design_3_Data.Range("E44").Value = 30.03126843510636
$ws_Design_3_Data.Range("H44").Value = 0.9081382662782423
$ws_Design_3_Data.Range("E45").Value = 29.29947468622886
$ws_Design_3_Data.Range("H45").Value = 7.208728303019477
$ws_Design_3_Data.Range("E46").Value = 8.45628332036219
$ws_Design_3_Data.Range("E47").Value = 507.9658430584618
$ws_Design_3_Data.Range("E48").Value = 63.7473665688842
$ws_Design_3_Data.Range("H48").Value = 31.87447604320883

# ---- Design 4 Data ----
$ws_Design_4_Data = $wb.Worksheets.Item("Design 4 Data")
$ws_Design_4_Data.Range("J13").Value = "max_fuel_L [L]"
$ws_Design_4_Data.Range("K13").Value = 57648.46529873692
$ws_Design_4_Data.Range("J14").Value = "total_fuel_L [L]"
$ws_Design_4_Data.Range("K14").Value = 52407.69572612447
$ws_Design_4_Data.Range("J15").Value = "mission_fuel_L [L]"
$ws_Design_4_Data.Range("K15").Value = 51903.21685673379
$ws_Design_4_Data.Range("J16").Value = "reserve_fuel_L [L]"
$ws_Design_4_Data.Range("K16").Value = 504.4788693906718
$ws_Design_4_Data.Range("J17").Value = "MTOM [kg]"
$ws_Design_4_Data.Range("K17").Value = 263578.3926925393
$ws_Design_4_Data.Range("J18").Value = "S [m²]"
$ws_Design_4_Data.Range("K18").Value = 561.0633135905323
$ws_Design_4_Data.Range("J19").Value = "b [m]"
$ws_Design_4_Data.Range("K19").Value = 41.02669790236105
$ws_Design_4_Data.Range("J20").Value = "MAC [m]"
$ws_Design_4_Data.Range("K20").Value = 13.67556596745369
$ws_Design_4_Data.Range("J21").Value = "fuel_economy [L/ton/km]"
$ws_Design_4_Data.Range("K21").Value = 0.07594985111932576
$ws_Design_4_Data.Range("J22").Value = "OEW [N]"
$ws_Design_4_Data.Range("K22").Value = 1274471.092321406
$ws_Design_4_Data.Range("J23").Value = "total_fuel [N]"
$ws_Design_4_Data.Range("K23").Value = 421577.9859600904
$ws_Design_4_Data.Range("J24").Value = "max_fuel [N]"
$ws_Design_4_Data.Range("K24").Value = 463735.7845560995
$ws_Design_4_Data.Range("J25").Value = "mission_fuel [N]"
$ws_Design_4_Data.Range("K25").Value = 417519.857038938
$ws_Design_4_Data.Range("J26").Value = "reserve_fuel [N]"
$ws_Design_4_Data.Range("K26").Value = 4058.128921152442
$ws_Design_4_Data.Range("K26").NumberFormat = "#,##0.000"
$ws_Design_4_Data.Range("J27").Value = "Mff"
$ws_Design_4_Data.Range("K27").Value = 0.8239224430163152
$ws_Design_4_Data.Range("K27").NumberFormat = "#,##0.000"
$ws_Design_4_Data.Range("J28").Value = "MTOW [N]"
$ws_Design_4_Data.Range("K28").Value = 2585704.03231381
$ws_Design_4_Data.Range("K28").NumberFormat = "#,##0.000"
$ws_Design_4_Data.Range("J29").Value = "LD"
$ws_Design_4_Data.Range("K29").Value = 23.95928996628532
$ws_Design_4_Data.Range("K29").NumberFormat = "#,##0.000"
$ws_Design_4_Data.Range("J30").Value = "total_fuel_L [L]"
$ws_Design_4_Data.Range("K30").Value = 52407.69572612447
$ws_Design_4_Data.Range("K30").NumberFormat = "#,##0.000"
$ws_Design_4_Data.Range("J31").Value = "mission_fuel_L [L]"
$ws_Design_4_Data.Range("K31").Value = 51903.21685673379
$ws_Design_4_Data.Range("K31").NumberFormat = "#,##0.000"
$ws_Design_4_Data.Range("J32").Value = "reserve_fuel_L [L]"
$ws_Design_4_Data.Range("K32").Value = 504.4788693906718
$ws_Design_4_Data.Range("K32").NumberFormat = "#,##0.000"
$ws_Design_4_Data.Range("J33").Value = "max_fuel_L [L]"
$ws_Design_4_Data.Range("K33").Value = 57648.46529873692
$ws_Design_4_Data.Range("K33").NumberFormat = "#,##0.000"
